$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sandeep konduri - hyperlinked email
$ws.Range("A2").Value = "Sandeep"
$ws.Range("B2").Value = "konduri"
$ws.Range("C2").Value = "sandeep.konduri123@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:sandeep.konduri123@gmail.com")

# Row 3: santhosh konduri - hyperlinked email
$ws.Range("A3").Value = "santhosh"
$ws.Range("B3").Value = "konduri"
$ws.Range("C3").Value = "konduri.santhosh96@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:konduri.santhosh96@gmail.com")

# Row 4: goutham - email styled like a hyperlink but no live link
$ws.Range("A4").Value = "goutham"
$ws.Range("C4").Value = "goutham221997@gmail.com"
$ws.Range("C4").Style = "Hyperlink"

# Row 5: only an email, styled like a hyperlink but no live link
$ws.Range("C5").Value = "vamsi.gch@gmail.com"
$ws.Range("C5").Style = "Hyperlink"

# Widen column C to fit the email addresses
$ws.Columns("C").ColumnWidth = 27.666666666666668

# Leave the final selection on A5, as in the saved workbook
$ws.Range("A5").Select() | Out-Null
